# Commit: "Added sssec for CD utility model."
#
# The data correction: the "q_dot_s_orig" value (column P) for the two
# CD-utility-model rows (row 3 = "Car - Germany", row 5 = "Lamp -
# Germany") was changed from 14425 to 12416. The author's on-screen
# selection/scroll position was also left on cell Q5 when the workbook
# was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EEU data")

# --- Data edits -------------------------------------------------------
$ws.Range("P3").Value = 12416
$ws.Range("P5").Value = 12416

# --- View state ---------------------------------------------------------
# Make sure the sheet is active, then move the selection to match where
# the author left the cursor when the file was saved.
$ws.Activate() | Out-Null
$ws.Range("Q5").Select() | Out-Null

# Best-effort: scroll the window so column B becomes the left-most
# visible column (topLeftCell="B1" in the saved view).
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
